$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The bullet "Experience on Perimeter security ... Paloalto, Cisco ASA,
# Juniper SRX, quagga, vyatta, VyOS etc ." drops "quagga, vyatta, VyOS"
# (the firewall/OS vendor is no longer listed) so it reads
# "... Paloalto, Cisco ASA, Juniper SRX etc .".
# ---------------------------------------------------------------------------

# 1) Remove "quagga, vyatta, VyOS" together with the following run's single
#    space, replacing with just one space (collapses several runs down to
#    the run that used to hold that trailing space).
$r1 = $d.Content
$found1 = $r1.Find.Execute("quagga, vyatta, VyOS ", $true, $false, $false, $false, $false, $true, 1, $false, " ", 2)
if (-not $found1) { throw "could not find 'quagga, vyatta, VyOS '" }

# 2) Drop the now-superfluous trailing ", " that used to introduce "quagga".
$r2 = $d.Content
$found2 = $r2.Find.Execute(", Cisco ASA, Juniper SRX, ", $true, $false, $false, $false, $false, $true, 1, $false, ", Cisco ASA, Juniper SRX", 2)
if (-not $found2) { throw "could not find ', Cisco ASA, Juniper SRX, '" }

# 3) " etc " is now contiguous text but the old spell-check markers
#    (proofErr spellStart/spellEnd) around "etc" are still sitting in the
#    run tree as stray siblings. Round-trip the text through a throwaway
#    value so the engine really rewrites this span (a set to the identical
#    string is treated as a no-op) which sweeps up the orphaned markers.
$r3 = $d.Content
$found3 = $r3.Find.Execute(" etc ", $true, $false, $false, $false, $false, $true, 1, $false, "##TEMP##", 2)
if (-not $found3) { throw "could not find ' etc '" }

$r4 = $d.Content
$found4 = $r4.Find.Execute("##TEMP##", $true, $false, $false, $false, $false, $true, 1, $false, " etc ", 2)
if (-not $found4) { throw "could not restore ' etc '" }

# 4) Re-split " etc " into its three original-looking runs and restore the
#    bold/complex-script-bold formatting that Word carries on this span:
#    " " (bold+bCs), "etc" (bCs only, not bold), " " (bold+bCs).
$s1 = $d.Content
$f1 = $s1.Find.Execute("Juniper SRX")
if (-not $f1) { throw "could not find 'Juniper SRX' anchor" }
$s1.Collapse(0)
$s1.MoveEnd(1, 1)
$s1.Font.BoldBi = 1

$s2 = $d.Content
$f2 = $s2.Find.Execute("Juniper SRX")
if (-not $f2) { throw "could not find 'Juniper SRX' anchor" }
$s2.Collapse(0)
$s2.MoveStart(1, 1)
$s2.MoveEnd(1, 3)
$s2.Font.Bold = 0
$s2.Font.BoldBi = 1

$s3 = $d.Content
$f3 = $s3.Find.Execute("Juniper SRX")
if (-not $f3) { throw "could not find 'Juniper SRX' anchor" }
$s3.Collapse(0)
$s3.MoveStart(1, 4)
$s3.MoveEnd(1, 1)
$s3.Font.BoldBi = 1

Write-Host "Done: removed 'quagga, vyatta, VyOS' from the firewall bullet."
